# Scheduled data refresh: update crafting-leve market price/profit figures
# (columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4022.25
$ws.Range("I76").Value = 3962.5925
$ws.Range("J76").Value = 4146.154
$ws.Range("K76").Value = 3962.5925
$ws.Range("L76").Value = 4146.154
$ws.Range("M76").Value = -3647.5925
$ws.Range("N76").Value = -4776.154
$ws.Range("H79").Value = 4022.25
$ws.Range("I79").Value = 3962.5925
$ws.Range("J79").Value = 4146.154
$ws.Range("K79").Value = 3962.5925
$ws.Range("L79").Value = 4146.154
$ws.Range("M79").Value = -2870.5925
$ws.Range("N79").Value = -6330.154
$ws.Range("H98").Value = 1636.409
$ws.Range("I98").Value = 1349.5555
$ws.Range("J98").Value = 2927.25
$ws.Range("K98").Value = 1349.5555
$ws.Range("L98").Value = 2927.25
$ws.Range("M98").Value = 148.4445000000001
$ws.Range("N98").Value = -5923.25
$ws.Range("H122").Value = 1636.409
$ws.Range("I122").Value = 1349.5555
$ws.Range("J122").Value = 2927.25
$ws.Range("K122").Value = 4048.6665
$ws.Range("L122").Value = 8781.75
$ws.Range("M122").Value = -1598.6665
$ws.Range("N122").Value = -13681.75
$ws.Range("H129").Value = 1001.6739
$ws.Range("J129").Value = 1064.5952
$ws.Range("L129").Value = 3193.7856
$ws.Range("N129").Value = -13193.7856
$ws.Range("H132").Value = 2128.3333
$ws.Range("I132").Value = 2128.3333
$ws.Range("K132").Value = 6384.999899999999
$ws.Range("M132").Value = -3854.999899999999
$ws.Range("H137").Value = 2440.7856
$ws.Range("I137").Value = 2634.375
$ws.Range("J137").Value = 2182.6667
$ws.Range("K137").Value = 7903.125
$ws.Range("L137").Value = 6548.000100000001
$ws.Range("M137").Value = -5353.125
$ws.Range("N137").Value = -11648.0001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3295.3333
$ws.Range("I63").Value = 2727.111
$ws.Range("K63").Value = 2727.111
$ws.Range("M63").Value = -2041.111
$ws.Range("H66").Value = 3295.3333
$ws.Range("I66").Value = 2727.111
$ws.Range("K66").Value = 13635.555
$ws.Range("M66").Value = -10203.555
$ws.Range("H88").Value = 14726.5
$ws.Range("I88").Value = 17968.666
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 17968.666
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -17562.666
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 14726.5
$ws.Range("I91").Value = 17968.666
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 17968.666
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -16564.666
$ws.Range("N91").Value = -7808
$ws.Range("H122").Value = 6946289.5
$ws.Range("I122").Value = 1774.9166
$ws.Range("J122").Value = 20835320
$ws.Range("K122").Value = 5324.7498
$ws.Range("L122").Value = 62505960
$ws.Range("M122").Value = -2874.7498
$ws.Range("N122").Value = -62510860
$ws.Range("H132").Value = 1962.9767
$ws.Range("I132").Value = 1777.5483
$ws.Range("J132").Value = 2442
$ws.Range("K132").Value = 5332.644899999999
$ws.Range("L132").Value = 7326
$ws.Range("M132").Value = -2802.644899999999
$ws.Range("N132").Value = -12386
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1878.3334
$ws.Range("I86").Value = 1864.5641
$ws.Range("K86").Value = 1864.5641
$ws.Range("M86").Value = -741.5641000000001
$ws.Range("H89").Value = 1878.3334
$ws.Range("I89").Value = 1864.5641
$ws.Range("K89").Value = 9322.8205
$ws.Range("M89").Value = -3706.8205
$ws.Range("H105").Value = 3658.2666
$ws.Range("I105").Value = 2844.0527
$ws.Range("J105").Value = 5064.636
$ws.Range("K105").Value = 2844.0527
$ws.Range("L105").Value = 5064.636
$ws.Range("M105").Value = -1097.0527
$ws.Range("N105").Value = -8558.636
$ws.Range("H134").Value = 2641.0977
$ws.Range("I134").Value = 2766.2424
$ws.Range("J134").Value = 2124.875
$ws.Range("K134").Value = 8298.727200000001
$ws.Range("L134").Value = 6374.625
$ws.Range("M134").Value = -5763.727200000001
$ws.Range("N134").Value = -11444.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 10561
$ws.Range("I122").Value = 6236.25
$ws.Range("K122").Value = 18708.75
$ws.Range("M122").Value = -16258.75
$ws.Range("H132").Value = 2660.1462
$ws.Range("I132").Value = 2376.074
$ws.Range("J132").Value = 3208
$ws.Range("K132").Value = 7128.222
$ws.Range("L132").Value = 9624
$ws.Range("M132").Value = -4598.222
$ws.Range("N132").Value = -14684
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 510
$ws.Range("I26").Value = 92.75
$ws.Range("J26").Value = 718.625
$ws.Range("K26").Value = 278.25
$ws.Range("L26").Value = 2155.875
$ws.Range("M26").Value = 9.75
$ws.Range("N26").Value = -2731.875
$ws.Range("H86").Value = 2416.6667
$ws.Range("I86").Value = 3234
$ws.Range("J86").Value = 1599.3334
$ws.Range("K86").Value = 9702
$ws.Range("L86").Value = 4798.0002
$ws.Range("M86").Value = -8516
$ws.Range("N86").Value = -7170.0002
$ws.Range("H87").Value = 8243.352999999999
$ws.Range("I87").Value = 999.6667
$ws.Range("K87").Value = 2999.0001
$ws.Range("M87").Value = -1751.0001
$ws.Range("H89").Value = 2416.6667
$ws.Range("I89").Value = 3234
$ws.Range("J89").Value = 1599.3334
$ws.Range("K89").Value = 29106
$ws.Range("L89").Value = 14394.0006
$ws.Range("M89").Value = -23178
$ws.Range("N89").Value = -26250.0006
$ws.Range("H90").Value = 8243.352999999999
$ws.Range("I90").Value = 999.6667
$ws.Range("K90").Value = 8997.0003
$ws.Range("M90").Value = -2757.0003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5569.8667
$ws.Range("I70").Value = 5433.143
$ws.Range("J70").Value = 5888.8887
$ws.Range("K70").Value = 5433.143
$ws.Range("L70").Value = 5888.8887
$ws.Range("M70").Value = -5163.143
$ws.Range("N70").Value = -6428.8887
$ws.Range("H73").Value = 5569.8667
$ws.Range("I73").Value = 5433.143
$ws.Range("J73").Value = 5888.8887
$ws.Range("K73").Value = 5433.143
$ws.Range("L73").Value = 5888.8887
$ws.Range("M73").Value = -4497.143
$ws.Range("N73").Value = -7760.8887
$ws.Range("H80").Value = 11886.667
$ws.Range("J80").Value = 7711.4287
$ws.Range("L80").Value = 7711.4287
$ws.Range("N80").Value = -9707.4287
$ws.Range("H83").Value = 11886.667
$ws.Range("J83").Value = 7711.4287
$ws.Range("L83").Value = 38557.14350000001
$ws.Range("N83").Value = -48541.14350000001
$ws.Range("H102").Value = 3755.5789
$ws.Range("I102").Value = 3354.6553
$ws.Range("K102").Value = 3354.6553
$ws.Range("M102").Value = -1732.6553
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6103.686
$ws.Range("I122").Value = 5256.2974
$ws.Range("J122").Value = 8343.214
$ws.Range("K122").Value = 15768.8922
$ws.Range("L122").Value = 25029.642
$ws.Range("M122").Value = -13318.8922
$ws.Range("N122").Value = -29929.642
